$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 55,6
$data[0,0] = 'Aalborg University'
$data[0,1] = 31
$data[0,2] = 25
$data[0,3] = 80.59999999999999
$data[0,4] = 63.7
$data[0,5] = 90.8
$data[1,0] = 'Aalborg University Hospital'
$data[1,1] = 35
$data[1,2] = 30
$data[1,3] = 85.7
$data[1,4] = 70.59999999999999
$data[1,5] = 93.7
$data[2,0] = 'Aarhus University'
$data[2,1] = 191
$data[2,2] = 158
$data[2,3] = 82.7
$data[2,4] = 76.7
$data[2,5] = 87.40000000000001
$data[3,0] = 'Aarhus University Hospital'
$data[3,1] = 66
$data[3,2] = 51
$data[3,3] = 77.3
$data[3,4] = 65.8
$data[3,5] = 85.7
$data[4,0] = 'Akershus University Hospital'
$data[4,1] = 11
$data[4,2] = 8
$data[4,3] = 72.7
$data[4,4] = 43.4
$data[4,5] = 90.3
$data[5,0] = 'Bispebjerg and Frederiksberg Hospital'
$data[5,1] = 77
$data[5,2] = 65
$data[5,3] = 84.40000000000001
$data[5,4] = 74.7
$data[5,5] = 90.90000000000001
$data[6,0] = 'Copenhagen University Hospital'
$data[6,1] = 172
$data[6,2] = 133
$data[6,3] = 77.3
$data[6,4] = 70.5
$data[6,5] = 82.89999999999999
$data[7,0] = 'Danderyd Hospital'
$data[7,1] = 9
$data[7,2] = 7
$data[7,3] = 77.8
$data[7,4] = 45.3
$data[7,5] = 96.09999999999999
$data[8,0] = 'Gothenburg University'
$data[8,1] = 25
$data[8,2] = 22
$data[8,3] = 88
$data[8,4] = 70
$data[8,5] = 95.8
$data[9,0] = 'Haukeland university hospital'
$data[9,1] = 30
$data[9,2] = 15
$data[9,3] = 50
$data[9,4] = 33.2
$data[9,5] = 66.8
$data[10,0] = 'Helsinki University Hospital'
$data[10,1] = 48
$data[10,2] = 41
$data[10,3] = 85.40000000000001
$data[10,4] = 72.8
$data[10,5] = 92.80000000000001
$data[11,0] = 'Herlev and Gentofte Hospital'
$data[11,1] = 111
$data[11,2] = 83
$data[11,3] = 74.8
$data[11,4] = 66
$data[11,5] = 81.89999999999999
$data[12,0] = 'Holbæk Hospital'
$data[12,1] = 3
$data[12,2] = 3
$data[12,3] = 100
$data[12,4] = 43.9
$data[12,5] = 100
$data[13,0] = 'Hvidovre and Amager Hospital'
$data[13,1] = 51
$data[13,2] = 43
$data[13,3] = 84.3
$data[13,4] = 72
$data[13,5] = 91.8
$data[14,0] = 'Karolinska Institutet'
$data[14,1] = 166
$data[14,2] = 131
$data[14,3] = 78.90000000000001
$data[14,4] = 72.09999999999999
$data[14,5] = 84.39999999999999
$data[15,0] = 'Karolinska University Hospital'
$data[15,1] = 48
$data[15,2] = 36
$data[15,3] = 75
$data[15,4] = 61.2
$data[15,5] = 85.09999999999999
$data[16,0] = 'Kuopio University Hospital'
$data[16,1] = 17
$data[16,2] = 14
$data[16,3] = 82.40000000000001
$data[16,4] = 59
$data[16,5] = 93.8
$data[17,0] = 'Linkoeping University'
$data[17,1] = 39
$data[17,2] = 31
$data[17,3] = 79.5
$data[17,4] = 64.5
$data[17,5] = 89.2
$data[18,0] = 'Linkoeping University Hospital'
$data[18,1] = 18
$data[18,2] = 14
$data[18,3] = 77.8
$data[18,4] = 54.8
$data[18,5] = 91
$data[19,0] = 'Lund University'
$data[19,1] = 36
$data[19,2] = 23
$data[19,3] = 63.9
$data[19,4] = 47.59999999999999
$data[19,5] = 77.5
$data[20,0] = 'Mental health services in the Capital Region of Denmark'
$data[20,1] = 4
$data[20,2] = 4
$data[20,3] = 100
$data[20,4] = 51
$data[20,5] = 100
$data[21,0] = 'Næstved Hospital'
$data[21,1] = 5
$data[21,2] = 5
$data[21,3] = 100
$data[21,4] = 56.59999999999999
$data[21,5] = 100
$data[22,0] = 'Nordsjællands Hospital'
$data[22,1] = 10
$data[22,2] = 9
$data[22,3] = 90
$data[22,4] = 59.59999999999999
$data[22,5] = 99.5
$data[23,0] = 'Norwegian University of Science and Technology'
$data[23,1] = 76
$data[23,2] = 53
$data[23,3] = 69.7
$data[23,4] = 58.7
$data[23,5] = 78.90000000000001
$data[24,0] = 'Odense University Hospital'
$data[24,1] = 91
$data[24,2] = 80
$data[24,3] = 87.90000000000001
$data[24,4] = 79.60000000000001
$data[24,5] = 93.10000000000001
$data[25,0] = 'Örebro University'
$data[25,1] = 17
$data[25,2] = 12
$data[25,3] = 70.59999999999999
$data[25,4] = 46.9
$data[25,5] = 86.7
$data[26,0] = 'Örebro University Hospital'
$data[26,1] = 1
$data[26,2] = 1
$data[26,3] = 100
$data[26,4] = 5.1
$data[26,5] = 100
$data[27,0] = 'Oslo University Hospital'
$data[27,1] = 102
$data[27,2] = 80
$data[27,3] = 78.40000000000001
$data[27,4] = 69.5
$data[27,5] = 85.3
$data[28,0] = 'Oulu University Hospital'
$data[28,1] = 10
$data[28,2] = 6
$data[28,3] = 60
$data[28,4] = 31.3
$data[28,5] = 83.2
$data[29,0] = 'Sahlgrenska University Hospital'
$data[29,1] = 42
$data[29,2] = 32
$data[29,3] = 76.2
$data[29,4] = 61.5
$data[29,5] = 86.5
$data[30,0] = 'Skane University Hospital'
$data[30,1] = 23
$data[30,2] = 19
$data[30,3] = 82.59999999999999
$data[30,4] = 62.9
$data[30,5] = 93
$data[31,0] = 'St. Olav’s University Hospital'
$data[31,1] = 24
$data[31,2] = 14
$data[31,3] = 58.3
$data[31,4] = 38.8
$data[31,5] = 75.5
$data[32,0] = 'Steno Diabetes Center Copenhagen'
$data[32,1] = 13
$data[32,2] = 11
$data[32,3] = 84.59999999999999
$data[32,4] = 57.8
$data[32,5] = 97.3
$data[33,0] = 'Stockholm South General Hospital'
$data[33,1] = 3
$data[33,2] = 3
$data[33,3] = 100
$data[33,4] = 43.9
$data[33,5] = 100
$data[34,0] = 'Tampere University Hospital'
$data[34,1] = 24
$data[34,2] = 20
$data[34,3] = 83.3
$data[34,4] = 64.09999999999999
$data[34,5] = 93.30000000000001
$data[35,0] = 'The National University Hospital of Iceland'
$data[35,1] = 5
$data[35,2] = 4
$data[35,3] = 80
$data[35,4] = 37.6
$data[35,5] = 99
$data[36,0] = 'Turku University Hospital'
$data[36,1] = 50
$data[36,2] = 30
$data[36,3] = 60
$data[36,4] = 46.2
$data[36,5] = 72.39999999999999
$data[37,0] = 'UiT The Arctic University of Norway'
$data[37,1] = 14
$data[37,2] = 10
$data[37,3] = 71.40000000000001
$data[37,4] = 45.4
$data[37,5] = 88.3
$data[38,0] = 'Umeå University'
$data[38,1] = 42
$data[38,2] = 29
$data[38,3] = 69
$data[38,4] = 54
$data[38,5] = 80.90000000000001
$data[39,0] = 'University Hospital of North Norway'
$data[39,1] = 17
$data[39,2] = 11
$data[39,3] = 64.7
$data[39,4] = 41.3
$data[39,5] = 82.69999999999999
$data[40,0] = 'University Hospital of Umeå'
$data[40,1] = 2
$data[40,2] = 2
$data[40,3] = 100
$data[40,4] = 17.8
$data[40,5] = 100
$data[41,0] = 'University of Bergen'
$data[41,1] = 31
$data[41,2] = 22
$data[41,3] = 71
$data[41,4] = 53.40000000000001
$data[41,5] = 83.89999999999999
$data[42,0] = 'University of Copenhagen'
$data[42,1] = 99
$data[42,2] = 73
$data[42,3] = 73.7
$data[42,4] = 64.3
$data[42,5] = 81.39999999999999
$data[43,0] = 'University of Eastern Finland'
$data[43,1] = 12
$data[43,2] = 8
$data[43,3] = 66.7
$data[43,4] = 39.1
$data[43,5] = 86.2
$data[44,0] = 'University of Helsinki'
$data[44,1] = 21
$data[44,2] = 17
$data[44,3] = 81
$data[44,4] = 60
$data[44,5] = 92.30000000000001
$data[45,0] = 'University of Iceland'
$data[45,1] = 5
$data[45,2] = 4
$data[45,3] = 80
$data[45,4] = 37.6
$data[45,5] = 99
$data[46,0] = 'University of Oslo'
$data[46,1] = 23
$data[46,2] = 21
$data[46,3] = 91.3
$data[46,4] = 73.2
$data[46,5] = 98.5
$data[47,0] = 'University of Oulu'
$data[47,1] = 25
$data[47,2] = 17
$data[47,3] = 68
$data[47,4] = 48.4
$data[47,5] = 82.8
$data[48,0] = 'University of Southern Denmark'
$data[48,1] = 42
$data[48,2] = 38
$data[48,3] = 90.5
$data[48,4] = 77.90000000000001
$data[48,5] = 96.2
$data[49,0] = 'University of Tampere'
$data[49,1] = 9
$data[49,2] = 9
$data[49,3] = 100
$data[49,4] = 70.09999999999999
$data[49,5] = 100
$data[50,0] = 'University of Turku'
$data[50,1] = 20
$data[50,2] = 15
$data[50,3] = 75
$data[50,4] = 53.1
$data[50,5] = 88.8
$data[51,0] = 'Uppsala Academic Hospital'
$data[51,1] = 9
$data[51,2] = 7
$data[51,3] = 77.8
$data[51,4] = 45.3
$data[51,5] = 96.09999999999999
$data[52,0] = 'Uppsala University'
$data[52,1] = 51
$data[52,2] = 33
$data[52,3] = 64.7
$data[52,4] = 51
$data[52,5] = 76.40000000000001
$data[53,0] = 'Zealand University Hospital'
$data[53,1] = 28
$data[53,2] = 24
$data[53,3] = 85.7
$data[53,4] = 68.5
$data[53,5] = 94.3
$data[54,0] = 'Total'
$data[54,1] = 2112
$data[54,2] = 1638
$data[54,3] = 77.60000000000001
$data[54,4] = 75.7
$data[54,5] = 79.3

$ws.Range("A2:F56").Value = $data

Write-Host "Done"